$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save current values for rows 17-20 (columns A-F) before rearranging.
$rowsToSave = 17, 18, 19, 20
$saved = @{}
foreach ($r in $rowsToSave) {
    $rowData = @{}
    foreach ($col in 1..6) {
        $cell = $ws.Cells.Item($r, $col)
        $rowData[$col] = @{
            Value = $cell.Value2
            Hyperlink = $null
        }
    }
    $saved[$r] = $rowData
}

# Save hyperlink info (address/sub-address/text) keyed by row, for column E.
$savedHyperlinks = @{}
foreach ($hl in $ws.Hyperlinks) {
    $hlRow = $hl.Range.Row
    if ($rowsToSave -contains $hlRow) {
        $savedHyperlinks[$hlRow] = @{
            Address = $hl.Address
            SubAddress = $hl.SubAddress
            TextToDisplay = $hl.TextToDisplay
        }
    }
}

# New row order: new row 17 gets old row 20's content, new 18 gets old 19,
# new 19 gets old 17, new 20 gets old 18.
$mapping = @{
    17 = 20
    18 = 19
    19 = 17
    20 = 18
}

# First, clear existing hyperlinks in these rows so we can re-add them cleanly.
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $hl = $ws.Hyperlinks.Item($i)
    $hlRow = $hl.Range.Row
    if ($rowsToSave -contains $hlRow) {
        $hl.Delete()
    }
}

# Write the values to their new rows.
foreach ($newRow in $rowsToSave) {
    $oldRow = $mapping[$newRow]
    $rowData = $saved[$oldRow]
    foreach ($col in 1..6) {
        $ws.Cells.Item($newRow, $col).Value2 = $rowData[$col].Value
    }
}

# Re-add hyperlinks (column E) to their new rows.
foreach ($newRow in $rowsToSave) {
    $oldRow = $mapping[$newRow]
    if ($savedHyperlinks.ContainsKey($oldRow)) {
        $hlInfo = $savedHyperlinks[$oldRow]
        $targetCell = $ws.Cells.Item($newRow, 5)
        $ws.Hyperlinks.Add($targetCell, $hlInfo.Address, $hlInfo.SubAddress, $null, $hlInfo.TextToDisplay) | Out-Null
    }
}

# Update the selection to match the target state: row 18 fully selected.
$ws.Range("A18:XFD18").Select()
